$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1936.7433  # H40: 1979.25 -> 1936.7433
$ws.Cells.Item(40, 9).Value = 1953.0483  # I40: 1984.7192 -> 1953.0483
$ws.Cells.Item(40, 10).Value = 1852.5  # J40: 1950.909 -> 1852.5
$ws.Cells.Item(40, 11).Value = 1953.0483  # K40: 1984.7192 -> 1953.0483
$ws.Cells.Item(40, 12).Value = 1852.5  # L40: 1950.909 -> 1852.5
$ws.Cells.Item(40, 13).Value = -1778.0483  # M40: -1809.7192 -> -1778.0483
$ws.Cells.Item(40, 14).Value = -2202.5  # N40: -2300.909 -> -2202.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 46296924  # H92: 53419570 -> 46296924
$ws.Cells.Item(92, 9).Value = 2137323.5  # I92: 2525971.8 -> 2137323.5
$ws.Cells.Item(92, 11).Value = 2137323.5  # K92: 2525971.8 -> 2137323.5
$ws.Cells.Item(92, 13).Value = -2136075.5  # M92: -2524723.8 -> -2136075.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(99, 8).Value = 722.8333  # H99: 765.5625 -> 722.8333
$ws.Cells.Item(99, 9).Value = 661.3333  # I99: 686.9091 -> 661.3333
$ws.Cells.Item(99, 10).Value = 845.8333  # J99: 938.6 -> 845.8333
$ws.Cells.Item(99, 11).Value = 1983.9999  # K99: 2060.7273 -> 1983.9999
$ws.Cells.Item(99, 12).Value = 2537.4999  # L99: 2815.8 -> 2537.4999
$ws.Cells.Item(99, 13).Value = -485.9999  # M99: -562.7273 -> -485.9999
$ws.Cells.Item(99, 14).Value = -5533.4999  # N99: -5811.8 -> -5533.4999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 14746076  # H113: 15974715 -> 14746076
$ws.Cells.Item(113, 9).Value = 6669198.5  # I113: 7409953 -> 6669198.5
$ws.Cells.Item(113, 11).Value = 6669198.5  # K113: 7409953 -> 6669198.5
$ws.Cells.Item(113, 13).Value = -6665944.5  # M113: -7406699 -> -6665944.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 1457.1321  # H137: 1517 -> 1457.1321
$ws.Cells.Item(137, 9).Value = 1436.5  # I137: 1520.3529 -> 1436.5
$ws.Cells.Item(137, 11).Value = 4309.5  # K137: 4561.0587 -> 4309.5
$ws.Cells.Item(137, 13).Value = -1759.5  # M137: -2011.0587 -> -1759.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 13494.875  # H45: 15278.143 -> 13494.875
$ws.Cells.Item(45, 9).Value = 13494.875  # I45: 15278.143 -> 13494.875
$ws.Cells.Item(45, 11).Value = 13494.875  # K45: 15278.143 -> 13494.875
$ws.Cells.Item(45, 13).Value = -13117.875  # M45: -14901.143 -> -13117.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 5825.8965  # H61: 3890.1914 -> 5825.8965
$ws.Cells.Item(61, 9).Value = 6933.55  # I61: 4014.7104 -> 6933.55
$ws.Cells.Item(61, 11).Value = 6933.55  # K61: 4014.7104 -> 6933.55
$ws.Cells.Item(61, 13).Value = -6721.55  # M61: -3802.7104 -> -6721.55

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(97, 8).Value = 1185.579  # H97: 1084.8636 -> 1185.579
$ws.Cells.Item(97, 9).Value = 943.3333  # I97: 810.7143 -> 943.3333
$ws.Cells.Item(97, 10).Value = 1600.8572  # J97: 1564.625 -> 1600.8572
$ws.Cells.Item(97, 11).Value = 943.3333  # K97: 810.7143 -> 943.3333
$ws.Cells.Item(97, 12).Value = 1600.8572  # L97: 1564.625 -> 1600.8572
$ws.Cells.Item(97, 13).Value = -447.3333  # M97: -314.7143 -> -447.3333
$ws.Cells.Item(97, 14).Value = -2592.8572  # N97: -2556.625 -> -2592.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 5825.8965  # H136: 3890.1914 -> 5825.8965
$ws.Cells.Item(136, 9).Value = 6933.55  # I136: 4014.7104 -> 6933.55
$ws.Cells.Item(136, 11).Value = 20800.65  # K136: 12044.1312 -> 20800.65
$ws.Cells.Item(136, 13).Value = -18250.65  # M136: -9494.1312 -> -18250.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(5, 8).Value = 2980.6667  # H5: 3971 -> 2980.6667
$ws.Cells.Item(5, 9).Value = 2976  # I5: 4002 -> 2976
$ws.Cells.Item(5, 10).Value = 2990  # J5: 3940 -> 2990
$ws.Cells.Item(5, 11).Value = 2976  # K5: 4002 -> 2976
$ws.Cells.Item(5, 12).Value = 2990  # L5: 3940 -> 2990
$ws.Cells.Item(5, 13).Value = -2863  # M5: -3889 -> -2863
$ws.Cells.Item(5, 14).Value = -3216  # N5: -4166 -> -3216

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2146.5  # H94: 1002.975 -> 2146.5
$ws.Cells.Item(94, 9).Value = 1293  # I94: 530 -> 1293
$ws.Cells.Item(94, 10).Value = 3000  # J94: 2249.9092 -> 3000
$ws.Cells.Item(94, 11).Value = 1293  # K94: 530 -> 1293
$ws.Cells.Item(94, 12).Value = 3000  # L94: 2249.9092 -> 3000
$ws.Cells.Item(94, 13).Value = -842  # M94: -79 -> -842
$ws.Cells.Item(94, 14).Value = -3902  # N94: -3151.9092 -> -3902

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 83334480  # H99: 250000900 -> 83334480
$ws.Cells.Item(99, 9).Value = 166667650  # I99: 500000350 -> 166667650
$ws.Cells.Item(99, 10).Value = 1316.8334  # J99: 1445 -> 1316.8334
$ws.Cells.Item(99, 11).Value = 166667650  # K99: 500000350 -> 166667650
$ws.Cells.Item(99, 12).Value = 1316.8334  # L99: 1445 -> 1316.8334
$ws.Cells.Item(99, 13).Value = -166666152  # M99: -499998852 -> -166666152
$ws.Cells.Item(99, 14).Value = -4312.8334  # N99: -4441 -> -4312.8334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 10944.875  # H105: 11377.305 -> 10944.875
$ws.Cells.Item(105, 9).Value = 23223.223  # I105: 26001.25 -> 23223.223
$ws.Cells.Item(105, 11).Value = 23223.223  # K105: 26001.25 -> 23223.223
$ws.Cells.Item(105, 13).Value = -21476.223  # M105: -24254.25 -> -21476.223

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1891.7059  # H16: 2042.2142 -> 1891.7059
$ws.Cells.Item(16, 9).Value = 1628.8334  # I16: 1775.3334 -> 1628.8334
$ws.Cells.Item(16, 11).Value = 1628.8334  # K16: 1775.3334 -> 1628.8334
$ws.Cells.Item(16, 13).Value = -1341.8334  # M16: -1488.3334 -> -1341.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 14461.667  # H31: 4353.8247 -> 14461.667
$ws.Cells.Item(31, 9).Value = 0  # I31: 1746.7028 -> 0
$ws.Cells.Item(31, 10).Value = 14461.667  # J31: 9177 -> 14461.667
$ws.Cells.Item(31, 11).Value = 0  # K31: 1746.7028 -> 0
$ws.Cells.Item(31, 12).Value = 14461.667  # L31: 9177 -> 14461.667
$ws.Cells.Item(31, 13).ClearContents()  # M31: -1451.7028 -> (removed)
$ws.Cells.Item(31, 14).Value = -15051.667  # N31: -9767 -> -15051.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 14461.667  # H34: 4353.8247 -> 14461.667
$ws.Cells.Item(34, 9).Value = 0  # I34: 1746.7028 -> 0
$ws.Cells.Item(34, 10).Value = 14461.667  # J34: 9177 -> 14461.667
$ws.Cells.Item(34, 11).Value = 0  # K34: 1746.7028 -> 0
$ws.Cells.Item(34, 12).Value = 14461.667  # L34: 9177 -> 14461.667
$ws.Cells.Item(34, 13).ClearContents()  # M34: -1544.7028 -> (removed)
$ws.Cells.Item(34, 14).Value = -14865.667  # N34: -9581 -> -14865.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1721.4445  # H58: 1827.1212 -> 1721.4445
$ws.Cells.Item(58, 9).Value = 1257.1111  # I58: 1396.7333 -> 1257.1111
$ws.Cells.Item(58, 11).Value = 1257.1111  # K58: 1396.7333 -> 1257.1111
$ws.Cells.Item(58, 13).Value = -1054.1111  # M58: -1193.7333 -> -1054.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 4402.8  # H86: 2302.6667 -> 4402.8
$ws.Cells.Item(86, 9).Value = 4335.6665  # I86: 2000 -> 4335.6665
$ws.Cells.Item(86, 10).Value = 4503.5  # J86: 2454 -> 4503.5
$ws.Cells.Item(86, 11).Value = 4335.6665  # K86: 2000 -> 4335.6665
$ws.Cells.Item(86, 12).Value = 4503.5  # L86: 2454 -> 4503.5
$ws.Cells.Item(86, 13).Value = -3212.6665  # M86: -877 -> -3212.6665
$ws.Cells.Item(86, 14).Value = -6749.5  # N86: -4700 -> -6749.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(89, 8).Value = 4402.8  # H89: 2302.6667 -> 4402.8
$ws.Cells.Item(89, 9).Value = 4335.6665  # I89: 2000 -> 4335.6665
$ws.Cells.Item(89, 10).Value = 4503.5  # J89: 2454 -> 4503.5
$ws.Cells.Item(89, 11).Value = 21678.3325  # K89: 10000 -> 21678.3325
$ws.Cells.Item(89, 12).Value = 22517.5  # L89: 12270 -> 22517.5
$ws.Cells.Item(89, 13).Value = -16062.3325  # M89: -4384 -> -16062.3325
$ws.Cells.Item(89, 14).Value = -33749.5  # N89: -23502 -> -33749.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 12346580  # H107: 12821440 -> 12346580
$ws.Cells.Item(107, 9).Value = 16667308  # I107: 17544522 -> 16667308
$ws.Cells.Item(107, 11).Value = 16667308  # K107: 17544522 -> 16667308
$ws.Cells.Item(107, 13).Value = -16665388  # M107: -17542602 -> -16665388

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1891.7059  # H113: 2042.2142 -> 1891.7059
$ws.Cells.Item(113, 9).Value = 1628.8334  # I113: 1775.3334 -> 1628.8334
$ws.Cells.Item(113, 11).Value = 1628.8334  # K113: 1775.3334 -> 1628.8334
$ws.Cells.Item(113, 13).Value = 541.1666  # M113: 394.6666 -> 541.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2573.111  # H132: 2593.8823 -> 2573.111
$ws.Cells.Item(132, 9).Value = 2374.8572  # I132: 2880 -> 2374.8572
$ws.Cells.Item(132, 10).Value = 2699.2727  # J132: 2474.6667 -> 2699.2727
$ws.Cells.Item(132, 11).Value = 7124.571599999999  # K132: 8640 -> 7124.571599999999
$ws.Cells.Item(132, 12).Value = 8097.8181  # L132: 7424.000100000001 -> 8097.8181
$ws.Cells.Item(132, 13).Value = -4594.571599999999  # M132: -6110 -> -4594.571599999999
$ws.Cells.Item(132, 14).Value = -13157.8181  # N132: -12484.0001 -> -13157.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 3336.8235  # H134: 2953.861 -> 3336.8235
$ws.Cells.Item(134, 9).Value = 3352.7856  # I134: 3060.862 -> 3352.7856
$ws.Cells.Item(134, 10).Value = 3262.3333  # J134: 2510.5715 -> 3262.3333
$ws.Cells.Item(134, 11).Value = 10058.3568  # K134: 9182.585999999999 -> 10058.3568
$ws.Cells.Item(134, 12).Value = 9786.999899999999  # L134: 7531.7145 -> 9786.999899999999
$ws.Cells.Item(134, 13).Value = -7523.356800000001  # M134: -6647.585999999999 -> -7523.356800000001
$ws.Cells.Item(134, 14).Value = -14856.9999  # N134: -12601.7145 -> -14856.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 1721.4445  # H136: 1827.1212 -> 1721.4445
$ws.Cells.Item(136, 9).Value = 1257.1111  # I136: 1396.7333 -> 1257.1111
$ws.Cells.Item(136, 11).Value = 3771.3333  # K136: 4190.199900000001 -> 3771.3333
$ws.Cells.Item(136, 13).Value = -1221.3333  # M136: -1640.199900000001 -> -1221.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(59, 8).Value = 5266.6665  # H59: 4316.6665 -> 5266.6665
$ws.Cells.Item(59, 9).Value = 2000  # I59: 3000 -> 2000
$ws.Cells.Item(59, 10).Value = 6900  # J59: 6950 -> 6900
$ws.Cells.Item(59, 11).Value = 6000  # K59: 9000 -> 6000
$ws.Cells.Item(59, 12).Value = 20700  # L59: 20850 -> 20700
$ws.Cells.Item(59, 13).Value = -5460  # M59: -8460 -> -5460
$ws.Cells.Item(59, 14).Value = -21780  # N59: -21930 -> -21780

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 2190.6072  # H122: 2419.4 -> 2190.6072
$ws.Cells.Item(122, 10).Value = 2417.6736  # J122: 2715.3953 -> 2417.6736
$ws.Cells.Item(122, 12).Value = 21759.0624  # L122: 24438.5577 -> 21759.0624
$ws.Cells.Item(122, 14).Value = -26659.0624  # N122: -29338.5577 -> -26659.0624

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 51391096  # H122: 51391124 -> 51391096
$ws.Cells.Item(122, 9).Value = 79061170  # I122: 79061220 -> 79061170
$ws.Cells.Item(122, 11).Value = 237183510  # K122: 237183660 -> 237183510
$ws.Cells.Item(122, 13).Value = -237181060  # M122: -237181210 -> -237181060

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value = 18841.297  # H123: 21377.703 -> 18841.297
$ws.Cells.Item(123, 10).Value = 18841.297  # J123: 21377.703 -> 18841.297
$ws.Cells.Item(123, 12).Value = 18841.297  # L123: 21377.703 -> 18841.297
$ws.Cells.Item(123, 14).Value = -23741.297  # N123: -26277.703 -> -23741.297

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4096.8335  # H126: 4288.74 -> 4096.8335
$ws.Cells.Item(126, 9).Value = 7988.375  # I126: 10085.167 -> 7988.375
$ws.Cells.Item(126, 11).Value = 23965.125  # K126: 30255.501 -> 23965.125
$ws.Cells.Item(126, 13).Value = -21495.125  # M126: -27785.501 -> -21495.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1857.6  # H61: 2056.6 -> 1857.6
$ws.Cells.Item(61, 9).Value = 1897.4286  # I61: 2122 -> 1897.4286
$ws.Cells.Item(61, 10).Value = 1764.6666  # J61: 1795 -> 1764.6666
$ws.Cells.Item(61, 11).Value = 1897.4286  # K61: 2122 -> 1897.4286
$ws.Cells.Item(61, 12).Value = 1764.6666  # L61: 1795 -> 1764.6666
$ws.Cells.Item(61, 13).Value = -1695.4286  # M61: -1920 -> -1695.4286
$ws.Cells.Item(61, 14).Value = -2168.6666  # N61: -2199 -> -2168.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 997.9231  # H93: 852.2727 -> 997.9231
$ws.Cells.Item(93, 9).Value = 781.8570999999999  # I93: 627.1539 -> 781.8570999999999
$ws.Cells.Item(93, 10).Value = 1250  # J93: 1177.4445 -> 1250
$ws.Cells.Item(93, 11).Value = 781.8570999999999  # K93: 627.1539 -> 781.8570999999999
$ws.Cells.Item(93, 12).Value = 1250  # L93: 1177.4445 -> 1250
$ws.Cells.Item(93, 13).Value = 466.1429000000001  # M93: 620.8461 -> 466.1429000000001
$ws.Cells.Item(93, 14).Value = -3746  # N93: -3673.4445 -> -3746

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 1857.6  # H113: 2056.6 -> 1857.6
$ws.Cells.Item(113, 9).Value = 1897.4286  # I113: 2122 -> 1897.4286
$ws.Cells.Item(113, 10).Value = 1764.6666  # J113: 1795 -> 1764.6666
$ws.Cells.Item(113, 11).Value = 1897.4286  # K113: 2122 -> 1897.4286
$ws.Cells.Item(113, 12).Value = 1764.6666  # L113: 1795 -> 1764.6666
$ws.Cells.Item(113, 13).Value = 272.5714  # M113: 48 -> 272.5714
$ws.Cells.Item(113, 14).Value = -6104.6666  # N113: -6135 -> -6104.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 2629200.5  # H122: 2547097.2 -> 2629200.5
$ws.Cells.Item(122, 9).Value = 3403514.2  # I122: 3248895.2 -> 3403514.2
$ws.Cells.Item(122, 11).Value = 10210542.6  # K122: 9746685.600000001 -> 10210542.6
$ws.Cells.Item(122, 13).Value = -10208092.6  # M122: -9744235.600000001 -> -10208092.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 10946.482  # H136: 7845.143 -> 10946.482
$ws.Cells.Item(136, 9).Value = 9099.352999999999  # I136: 5642.6553 -> 9099.352999999999
$ws.Cells.Item(136, 10).Value = 13563.25  # J136: 12758.385 -> 13563.25
$ws.Cells.Item(136, 11).Value = 27298.059  # K136: 16927.9659 -> 27298.059
$ws.Cells.Item(136, 12).Value = 40689.75  # L136: 38275.155 -> 40689.75
$ws.Cells.Item(136, 13).Value = -24748.059  # M136: -14377.9659 -> -24748.059

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2333.1667  # H96: 2420 -> 2333.1667
$ws.Cells.Item(96, 9).Value = 2399.8  # I96: 2525 -> 2399.8
$ws.Cells.Item(96, 11).Value = 2399.8  # K96: 2525 -> 2399.8
$ws.Cells.Item(96, 13).Value = -1026.8  # M96: -1152 -> -1026.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1099.5  # H126: 1142.2858 -> 1099.5
$ws.Cells.Item(126, 9).Value = 807.6667  # I126: 809.2 -> 807.6667
$ws.Cells.Item(126, 11).Value = 2423.0001  # K126: 2427.6 -> 2423.0001
$ws.Cells.Item(126, 13).Value = 46.9998999999998  # M126: 42.39999999999964 -> 46.9998999999998
